# Auto-generated script applying scheduled market-data refresh to Mandragora_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 316.5
$ws.Range("I38").Value = 276
$ws.Range("J38").Value = 600
$ws.Range("K38").Value = 828
$ws.Range("L38").Value = 1800
$ws.Range("M38").Value = -456
$ws.Range("N38").Value = -2544
$ws.Range("H40").Value = 2417.8667
$ws.Range("I40").Value = 1750.2858
$ws.Range("K40").Value = 1750.2858
$ws.Range("M40").Value = -1575.2858
$ws.Range("H64").Value = 3799.6667
$ws.Range("I64").Value = 3521.4285
$ws.Range("J64").Value = 3938.7856
$ws.Range("K64").Value = 3521.4285
$ws.Range("L64").Value = 3938.7856
$ws.Range("M64").Value = -3273.4285
$ws.Range("N64").Value = -4434.7856
$ws.Range("H67").Value = 3799.6667
$ws.Range("I67").Value = 3521.4285
$ws.Range("J67").Value = 3938.7856
$ws.Range("K67").Value = 3521.4285
$ws.Range("L67").Value = 3938.7856
$ws.Range("M67").Value = -2663.4285
$ws.Range("N67").Value = -5654.7856
$ws.Range("H112").Value = 2841.4
$ws.Range("J112").Value = 3178.6333
$ws.Range("L112").Value = 9535.8999
$ws.Range("N112").Value = -11751.8999
$ws.Range("H113").Value = 2807.724
$ws.Range("I113").Value = 1842.9
$ws.Range("J113").Value = 4951.778
$ws.Range("K113").Value = 1842.9
$ws.Range("L113").Value = 4951.778
$ws.Range("M113").Value = 1411.1
$ws.Range("N113").Value = -11459.778
$ws.Range("H138").Value = 2466.1064
$ws.Range("I138").Value = 1179.75
$ws.Range("J138").Value = 4361.7896
$ws.Range("K138").Value = 3539.25
$ws.Range("L138").Value = 13085.3688
$ws.Range("M138").Value = 1600.75
$ws.Range("N138").Value = -23365.3688

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1520.2106
$ws.Range("I74").Value = 1134.2858
$ws.Range("J74").Value = 2600.8
$ws.Range("K74").Value = 1134.2858
$ws.Range("L74").Value = 2600.8
$ws.Range("M74").Value = -260.2858000000001
$ws.Range("N74").Value = -4348.8
$ws.Range("H77").Value = 1520.2106
$ws.Range("I77").Value = 1134.2858
$ws.Range("J77").Value = 2600.8
$ws.Range("K77").Value = 5671.429
$ws.Range("L77").Value = 13004
$ws.Range("M77").Value = -1303.429
$ws.Range("N77").Value = -21740
$ws.Range("H104").Value = 20225
$ws.Range("J104").Value = 20225
$ws.Range("L104").Value = 20225
$ws.Range("N104").Value = -27213

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 12348612
$ws.Range("I31").Value = 2170.8572
$ws.Range("J31").Value = 25644780
$ws.Range("K31").Value = 2170.8572
$ws.Range("L31").Value = 25644780
$ws.Range("M31").Value = -1875.8572
$ws.Range("N31").Value = -25645370
$ws.Range("H34").Value = 12348612
$ws.Range("I34").Value = 2170.8572
$ws.Range("J34").Value = 25644780
$ws.Range("K34").Value = 2170.8572
$ws.Range("L34").Value = 25644780
$ws.Range("M34").Value = -1968.8572
$ws.Range("N34").Value = -25645184

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 668.9091
$ws.Range("I5").Value = 247.8
$ws.Range("K5").Value = 743.4000000000001
$ws.Range("M5").Value = -631.4000000000001
$ws.Range("H131").Value = 1114.2609
$ws.Range("J131").Value = 1157.4
$ws.Range("L131").Value = 3472.2
$ws.Range("N131").Value = -13552.2
$ws.Range("H135").Value = 668.9091
$ws.Range("I135").Value = 247.8
$ws.Range("K135").Value = 2230.2
$ws.Range("M135").Value = 304.7999999999997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 3600.7778
$ws.Range("I3").Value = 502.33334
$ws.Range("J3").Value = 5150
$ws.Range("K3").Value = 502.33334
$ws.Range("L3").Value = 5150
$ws.Range("M3").Value = -386.33334
$ws.Range("N3").Value = -5382
$ws.Range("H6").Value = 3490
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 3490
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 3490
$ws.Range("M6").ClearContents()
$ws.Range("N6").Value = -3716
$ws.Range("H7").Value = 5004000
$ws.Range("J7").Value = 5004000
$ws.Range("L7").Value = 5004000
$ws.Range("N7").Value = -5004224
$ws.Range("H8").Value = 5004000
$ws.Range("J8").Value = 5004000
$ws.Range("L8").Value = 5004000
$ws.Range("N8").Value = -5004278
$ws.Range("H9").Value = 1000
$ws.Range("I9").Value = 1000
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 1000
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = -830
$ws.Range("N9").ClearContents()
$ws.Range("H13").Value = 349.85715
$ws.Range("I13").Value = 349.85715
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 349.85715
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -210.85715
$ws.Range("N13").ClearContents()
$ws.Range("H14").Value = 116333570
$ws.Range("I14").Value = 174500000
$ws.Range("J14").Value = 700
$ws.Range("K14").Value = 174500000
$ws.Range("L14").Value = 700
$ws.Range("M14").Value = -174499832
$ws.Range("N14").Value = -1036
$ws.Range("H16").Value = 3490
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 3490
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 3490
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -3990
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 924.8182
$ws.Range("I22").Value = 700.3333
$ws.Range("J22").Value = 1009
$ws.Range("K22").Value = 700.3333
$ws.Range("L22").Value = 1009
$ws.Range("M22").Value = -405.3333
$ws.Range("N22").Value = -1599
$ws.Range("H27").Value = 924.8182
$ws.Range("I27").Value = 700.3333
$ws.Range("J27").Value = 1009
$ws.Range("K27").Value = 700.3333
$ws.Range("L27").Value = 1009
$ws.Range("M27").Value = -593.3333
$ws.Range("N27").Value = -1223
$ws.Range("H55").Value = 331.75
$ws.Range("I55").Value = 269
$ws.Range("J55").Value = 520
$ws.Range("K55").Value = 269
$ws.Range("L55").Value = 520
$ws.Range("M55").Value = -96
$ws.Range("N55").Value = -866
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H100").Value = 3490.5715
$ws.Range("J100").Value = 3700.1428
$ws.Range("L100").Value = 3700.1428
$ws.Range("N100").Value = -4782.1428

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 13685478
$ws.Range("J5").Value = 15396125
$ws.Range("L5").Value = 15396125
$ws.Range("N5").Value = -15396349
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("M19").ClearContents()
$ws.Range("N19").ClearContents()
$ws.Range("H80").Value = 19999
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 19999
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 19999
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -21995
$ws.Range("H83").Value = 19999
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 19999
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 59997
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -69981
$ws.Range("H107").Value = 443
$ws.Range("I107").Value = 402.2143
$ws.Range("J107").Value = 633.3333
$ws.Range("K107").Value = 1206.6429
$ws.Range("L107").Value = 1899.9999
$ws.Range("M107").Value = 713.3571000000002
$ws.Range("N107").Value = -5739.9999
$ws.Range("H109").Value = 25377
$ws.Range("J109").Value = 25377
$ws.Range("L109").Value = 25377
$ws.Range("N109").Value = -28151
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()

Write-Output "Applied 214 cell updates across 7 worksheets"
